$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos - updated stats ---
$ws.Cells.Item(4, 2).Value = 677056   # B4 Casos totales
$ws.Cells.Item(4, 3).Value = 29053    # C4 Nuevos casos
$ws.Cells.Item(4, 5).Value = 585205   # E4 Recuperados
$ws.Cells.Item(4, 7).Value = 2137     # G4 Casos criticos
$ws.Cells.Item(4, 8).Value = 34580    # H4 Muertes

# --- Row 14: Brasil - updated stats ---
$ws.Cells.Item(14, 2).Value = 30683   # B14 Casos totales
$ws.Cells.Item(14, 3).Value = 2073    # C14 Nuevos casos
$ws.Cells.Item(14, 5).Value = 14710   # E14 Recuperados
$ws.Cells.Item(14, 7).Value = 190     # G14 Casos criticos
$ws.Cells.Item(14, 8).Value = 1947    # H14 Muertes

# --- Rows 36-37: Chequia moves ahead of Arabia Saudita, both get refreshed values ---
$ws.Cells.Item(36, 1).Value = "Chequia"
$ws.Cells.Item(36, 2).Value = 6433
$ws.Cells.Item(36, 3).Value = 132
$ws.Cells.Item(36, 4).Value = 972
$ws.Cells.Item(36, 5).Value = 5292
$ws.Cells.Item(36, 6).Value = 75
$ws.Cells.Item(36, 7).Value = 3
$ws.Cells.Item(36, 8).Value = 169

$ws.Cells.Item(37, 1).Value = "Arabia Saudita"
$ws.Cells.Item(37, 2).Value = 6380
$ws.Cells.Item(37, 3).Value = 518
$ws.Cells.Item(37, 4).Value = 990
$ws.Cells.Item(37, 5).Value = 5307
$ws.Cells.Item(37, 6).Value = 71
$ws.Cells.Item(37, 7).Value = 4
$ws.Cells.Item(37, 8).Value = 83

# --- Row 120: Venezuela - updated stats ---
$ws.Cells.Item(120, 2).Value = 204    # B120 Casos totales
$ws.Cells.Item(120, 3).Value = 7      # C120 Nuevos casos
$ws.Cells.Item(120, 5).Value = 84     # E120 Recuperados
$ws.Cells.Item(120, 6).Value = 4      # F120 Casos criticos

# --- Row 157: Libia - updated stats ---
$ws.Cells.Item(157, 2).Value = 49     # B157 Casos totales
$ws.Cells.Item(157, 3).Value = 1      # C157 Nuevos casos
$ws.Cells.Item(157, 5).Value = 37     # E157 Recuperados
